# "orden de las carpetas de admin y uniones"
#
# Target changes for bd/pedidos.xlsx:
#   - Header row (A1:D1) loses the bold / bordered / centered style that was
#     applied to it (cells fall back to the workbook's default/"Normal"
#     style, index 0 - no more explicit font/border formatting).
#   - D2 ("estado" for the single order row) changes from "pagado" to
#     "enviado".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop all explicit formatting (bold font + thin border + centered
# alignment) from the header cells so they go back to the default style.
$ws.Range("A1:D1").ClearFormats()

# Update the order status value.
$ws.Range("D2").Value = "enviado"
